$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "43.221.64"
Set-TextValue $ws.Range("E2") "  -1.32%  "
Set-TextValue $ws.Range("D3") "2.272.71"
Set-TextValue $ws.Range("E3") "  -1.76%  "
Set-TextValue $ws.Range("E4") "  -0.03%  "
Set-TextValue $ws.Range("D5") "111.40"
Set-TextValue $ws.Range("E5") "  -1.20%  "
Set-TextValue $ws.Range("D6") "264.54"
Set-TextValue $ws.Range("E6") "  -2.27%  "
Set-TextValue $ws.Range("D7") "0.618"
Set-TextValue $ws.Range("E7") "  -1.56%  "
Set-TextValue $ws.Range("E8") "  +0.27%  "
Set-TextValue $ws.Range("D9") "0.604"
Set-TextValue $ws.Range("E9") "  -3.39%  "
Set-TextValue $ws.Range("D10") "47.54"
Set-TextValue $ws.Range("E10") "  +0.01%  "
Set-TextValue $ws.Range("D11") "0.0930"
Set-TextValue $ws.Range("E11") "  -1.51%  "
Set-TextValue $ws.Range("D12") "8.77"
Set-TextValue $ws.Range("E12") "  -1.22%  "
Set-TextValue $ws.Range("E13") "  +0.83%  "
Set-TextValue $ws.Range("D14") "15.41"
Set-TextValue $ws.Range("E14") "  -2.36%  "
Set-TextValue $ws.Range("D15") "2.608.49"
Set-TextValue $ws.Range("E15") "  -2.02%  "
Set-TextValue $ws.Range("D16") "0.851"
Set-TextValue $ws.Range("E16") "  -1.33%  "
Set-TextValue $ws.Range("D17") "2.264.01"
Set-TextValue $ws.Range("E17") "  -2.32%  "
Set-TextValue $ws.Range("D18") "43.058.27"
Set-TextValue $ws.Range("E18") "  -1.67%  "
Set-TextValue $ws.Range("D19") "0.0000108"
Set-TextValue $ws.Range("E19") "  -2.13%  "
Set-TextValue $ws.Range("D20") "6.81"
Set-TextValue $ws.Range("E20") "  +1.62%  "
Set-TextValue $ws.Range("D21") "71.06"
Set-TextValue $ws.Range("E21") "  -2.03%  "
Set-TextValue $ws.Range("D22") "2.51"
Set-TextValue $ws.Range("E22") "  +0.26%  "
Set-TextValue $ws.Range("D23") "229.83"
Set-TextValue $ws.Range("E23") "  -1.92%  "
Set-TextValue $ws.Range("D24") "9.65"
Set-TextValue $ws.Range("E24") "  +1.60%  "
Set-TextValue $ws.Range("D25") "2.86"
Set-TextValue $ws.Range("E25") "  -1.04%  "
Set-TextValue $ws.Range("E26") "  +0.07%  "
Set-TextValue $ws.Range("D27") "11.29"
Set-TextValue $ws.Range("E27") "  -2.04%  "
Set-TextValue $ws.Range("D28") "3.92"
Set-TextValue $ws.Range("E28") "  -1.10%  "
Set-TextValue $ws.Range("D29") "40.32"
Set-TextValue $ws.Range("E29") "  -5.14%  "
Set-TextValue $ws.Range("E30") "  -1.77%  "
Set-TextValue $ws.Range("E31") "  -4.77%  "
Set-TextValue $ws.Range("D32") "171.68"
Set-TextValue $ws.Range("E32") "  -3.28%  "
Set-TextValue $ws.Range("D33") "21.31"
Set-TextValue $ws.Range("E33") "  -3.08%  "
Set-TextValue $ws.Range("D34") "0.0902"
Set-TextValue $ws.Range("E34") "  -3.55%  "
Set-TextValue $ws.Range("D35") "5.71"
Set-TextValue $ws.Range("E35") "  +1.89%  "
Set-TextValue $ws.Range("E36") "  -0.44%  "
Set-TextValue $ws.Range("D37") "4.67"
Set-TextValue $ws.Range("E37") "  -2.28%  "
Set-TextValue $ws.Range("D38") "0.0351"
Set-TextValue $ws.Range("E38") "  -2.58%  "
Set-TextValue $ws.Range("E39") "  -6.39%  "
Set-TextValue $ws.Range("D40") "3.80"
Set-TextValue $ws.Range("E40") "  -4.21%  "
Set-TextValue $ws.Range("D41") "2.62"
Set-TextValue $ws.Range("E41") "  +9.19%  "
Set-TextValue $ws.Range("D42") "75.98"
Set-TextValue $ws.Range("E42") "  +9.34%  "
Set-TextValue $ws.Range("D43") "13.83"
Set-TextValue $ws.Range("E43") "  +7.40%  "
Set-TextValue $ws.Range("D44") "0.235"
Set-TextValue $ws.Range("E44") "  -4.03%  "
Set-TextValue $ws.Range("D45") "6.07"
Set-TextValue $ws.Range("E45") "  +4.40%  "
Set-TextValue $ws.Range("E46") "  +0.08%  "
Set-TextValue $ws.Range("D47") "1.37"
Set-TextValue $ws.Range("E47") "  -3.64%  "
Set-TextValue $ws.Range("D48") "8.62"
Set-TextValue $ws.Range("E48") "  -2.49%  "
Set-TextValue $ws.Range("D51") "100.87"
Set-TextValue $ws.Range("E51") "  +0.84%  "

# Rows 49/50 swap content: Cronos <-> TrustWalletToken
Set-TextValue $ws.Range("B49") "TrustWalletToken"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue $ws.Range("D49") "1.25"
Set-TextValue $ws.Range("E49") "  +1.87%  "

Set-TextValue $ws.Range("B50") "Cronos"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D50") "0.0991"
Set-TextValue $ws.Range("E50") "  -1.94%  "
